$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Cells.Item(2, 8).Value = 5754.1763
$v2 = $ws.Cells.Item(2, 8).Value2
Write-Host ("H2=" + $v2)
